$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7008573412895203
$ws.Range("B1").Value = 1.132298588752747
$ws.Range("C1").Value = 3.882263660430908
$ws.Range("D1").Value = 3.3177330493927
$ws.Range("E1").Value = 1.790637254714966
